# Update the offshore wind maximum-capacity values per cluster.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 4
$ws.Range("B4").Value = 53
$ws.Range("B5").Value = 6
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 0
$ws.Range("B8").Value = 18
$ws.Range("B9").Value = 0
